# daily auto push: 2025-10-09 18:38 UTC
# Append the new daily data row (row 86) to the bottom of Sheet1's table,
# extending the used range from A1:D85 to A1:D86.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 86

# Column A holds a date string such as "2025/09/22" stored as literal text
# (not an actual Excel date) throughout the sheet. Assigning a date-shaped
# string straight to .Value makes Excel auto-convert it to a date serial,
# which would diverge from every other row in the column. Instead, compute
# the text via TEXT()/"@" (forces a string result) and then paste-special
# just the value, exactly like using Excel's "Paste Values" after typing a
# formula - this keeps the stored cell a plain string, matching the rest of
# the column, with no extra number formatting applied to the cell.
$dateCell = $ws.Cells.Item($newRow, 1)
$dateCell.Formula = "=TEXT(""2025/10/10"", ""@"")"
$dateCell.Copy()
$dateCell.PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false

$ws.Cells.Item($newRow, 2).Value = "金"
$ws.Cells.Item($newRow, 3).Value = 1
$ws.Cells.Item($newRow, 4).Value = 201
